$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) value updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: move the value that was in D2 (39.05...) -> C2 gets the new measured value,
# D2 is cleared
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 39.663817176861357

# Row 3: B3 and C3 values removed
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Match the updated active selection on the sheet
$ws.Range("B1:E3").Select()
